$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B and C columns for rows 2-48 with newly computed values (improved performance via priority queue)
$ws.Cells.Item(2, 2).Value = 1.287483147441535
$ws.Cells.Item(2, 3).Value = 1.403005554198357
$ws.Cells.Item(3, 2).Value = 2.191553290604237
$ws.Cells.Item(3, 3).Value = 2.573896860326975
$ws.Cells.Item(4, 2).Value = 3.543790864104566
$ws.Cells.Item(4, 3).Value = 3.993392863159283
$ws.Cells.Item(5, 2).Value = 4.674828760453565
$ws.Cells.Item(5, 3).Value = 5.593223446315198
$ws.Cells.Item(6, 2).Value = 6.156779154688342
$ws.Cells.Item(6, 3).Value = 7.218194867730951
$ws.Cells.Item(7, 2).Value = 6.553875012213776
$ws.Cells.Item(7, 3).Value = 8.369258834213054
$ws.Cells.Item(8, 2).Value = 6.691653979187845
$ws.Cells.Item(8, 3).Value = 9.888812823761139
$ws.Cells.Item(9, 2).Value = 12.78044370847116
$ws.Cells.Item(9, 3).Value = 11.20048918982316
$ws.Cells.Item(10, 2).Value = 13.14934173963031
$ws.Cells.Item(10, 3).Value = 12.54584063506072
$ws.Cells.Item(11, 2).Value = 14.52428402938353
$ws.Cells.Item(11, 3).Value = 14.17422300672925
$ws.Cells.Item(12, 2).Value = 14.68089998200563
$ws.Cells.Item(12, 3).Value = 15.79405604972997
$ws.Cells.Item(13, 2).Value = 23.18683475578371
$ws.Cells.Item(13, 3).Value = 17.54096994485423
$ws.Cells.Item(14, 2).Value = 24.57264523443997
$ws.Cells.Item(14, 3).Value = 18.74422592060859
$ws.Cells.Item(15, 2).Value = 26.02710242037035
$ws.Cells.Item(15, 3).Value = 20.46841043523107
$ws.Cells.Item(16, 2).Value = 27.84525455872748
$ws.Cells.Item(16, 3).Value = 21.74844066490082
$ws.Cells.Item(17, 2).Value = 28.85769346028644
$ws.Cells.Item(17, 3).Value = 23.90869188078957
$ws.Cells.Item(18, 2).Value = 29.86862092005304
$ws.Cells.Item(18, 3).Value = 25.02280357817771
$ws.Cells.Item(19, 2).Value = 30.24721875574588
$ws.Cells.Item(19, 3).Value = 26.54228154854502
$ws.Cells.Item(20, 2).Value = 31.9796986071847
$ws.Cells.Item(20, 3).Value = 28.42884599299049
$ws.Cells.Item(21, 2).Value = 35.38388867374369
$ws.Cells.Item(21, 3).Value = 29.77841348580744
$ws.Cells.Item(22, 2).Value = 36.39557491615464
$ws.Cells.Item(22, 3).Value = 31.03479184445516
$ws.Cells.Item(23, 2).Value = 36.64729614751707
$ws.Cells.Item(23, 3).Value = 32.46723529721053
$ws.Cells.Item(24, 2).Value = 40.46893265650604
$ws.Cells.Item(24, 3).Value = 34.03649352120321
$ws.Cells.Item(25, 2).Value = 40.88419104277116
$ws.Cells.Item(25, 3).Value = 36.16930952749792
$ws.Cells.Item(26, 2).Value = 41.46108745763964
$ws.Cells.Item(26, 3).Value = 37.60719802707315
$ws.Cells.Item(27, 2).Value = 43.45171250100093
$ws.Cells.Item(27, 3).Value = 39.0379155722072
$ws.Cells.Item(28, 2).Value = 43.76161351304479
$ws.Cells.Item(28, 3).Value = 40.62186484454519
$ws.Cells.Item(29, 2).Value = 45.0625265750013
$ws.Cells.Item(29, 3).Value = 42.31454643332894
$ws.Cells.Item(30, 2).Value = 50.06535754112878
$ws.Cells.Item(30, 3).Value = 43.51303696249746
$ws.Cells.Item(31, 2).Value = 51.11547314149358
$ws.Cells.Item(31, 3).Value = 45.58383514304137
$ws.Cells.Item(32, 2).Value = 53.43446407192612
$ws.Cells.Item(32, 3).Value = 46.79461015563245
$ws.Cells.Item(33, 2).Value = 55.16186291924704
$ws.Cells.Item(33, 3).Value = 48.61013337034979
$ws.Cells.Item(34, 2).Value = 55.99721147052355
$ws.Cells.Item(34, 3).Value = 50.28261904336615
$ws.Cells.Item(35, 2).Value = 56.67181856806376
$ws.Cells.Item(35, 3).Value = 51.67553025662125
$ws.Cells.Item(36, 2).Value = 60.69452796372011
$ws.Cells.Item(36, 3).Value = 53.26608168250775
$ws.Cells.Item(37, 2).Value = 65.35007596665054
$ws.Cells.Item(37, 3).Value = 54.47323662478596
$ws.Cells.Item(38, 2).Value = 66.08632965479499
$ws.Cells.Item(38, 3).Value = 56.1350349999621
$ws.Cells.Item(39, 2).Value = 70.20159224790058
$ws.Cells.Item(39, 3).Value = 57.39973701877532
$ws.Cells.Item(40, 2).Value = 72.78342661002772
$ws.Cells.Item(40, 3).Value = 58.84689887551348
$ws.Cells.Item(41, 2).Value = 74.66460975615436
$ws.Cells.Item(41, 3).Value = 60.45322237527225
$ws.Cells.Item(42, 2).Value = 75.16025185232623
$ws.Cells.Item(42, 3).Value = 61.98918190035944
$ws.Cells.Item(43, 2).Value = 77.84837243418728
$ws.Cells.Item(43, 3).Value = 63.84285501948462
$ws.Cells.Item(44, 2).Value = 81.3302313536736
$ws.Cells.Item(44, 3).Value = 65.32829588116344
$ws.Cells.Item(45, 2).Value = 81.44140367820123
$ws.Cells.Item(45, 3).Value = 66.68751545914397
$ws.Cells.Item(46, 2).Value = 83.50592692063815
$ws.Cells.Item(46, 3).Value = 68.47478627072276
$ws.Cells.Item(47, 2).Value = 85.70882314310553
$ws.Cells.Item(47, 3).Value = 70.14026496508136
$ws.Cells.Item(48, 2).Value = 87.56062292395593
$ws.Cells.Item(48, 3).Value = 72.13287407628231

# Add new row 49 with the next simulation step
$ws.Cells.Item(49, 1).Value = 47
$ws.Cells.Item(49, 2).Value = 93.14206198308786
$ws.Cells.Item(49, 3).Value = 73.50679205353029

# Match the formatting (bold, bordered, centered style) used by the other rows in column A
$ws.Range("A48").Copy()
$ws.Range("A49").PasteSpecial(-4122)
$excel.CutCopyMode = 0
